$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q2:Q4 currently hold the merchant/card acceptor code as the text string
# "436845558641121"; replace each with the plain numeric value 1121.
$ws.Range("Q2").Value = 1121
$ws.Range("Q3").Value = 1121
$ws.Range("Q4").Value = 1121

# Update the active cell selection to match the authored state.
$ws.Range("Q3").Select()
